# Removed dependency of dose multiplier and updated documentation.
#
# The "Dose multiplier interval" column (column H, header shared-string
# index 9) is deleted entirely. Content shifts left: the old column I
# ("Force delay") becomes the new column H, and every column to the right
# shifts down by one letter (the sheet's used range shrinks from AF to AE).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole column H (shifts everything after it one column left).
$ws.Columns.Item(8).Delete()

# Match Excel's default post-delete selection: the column that now occupies
# the deleted column's position gets selected in full.
$ws.Columns.Item(8).Select()
